$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label in J11: "факт" -> "факт %$"
$ws.Range("J11").Value = "факт %$"

# Add new row 26 of data
# Copy the date format from F12 (style index 1) onto F26 before writing,
# so it reuses the existing style instead of creating a duplicate one.
$ws.Range("F12").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Value = 45623

$ws.Range("G26").Value = 3379.7

$ws.Range("H26").Formula = "=G26-`$D`$12"

$ws.Range("I26").Formula = "=H26/`$D`$14"

$ws.Range("J26").Value = 23

$ws.Range("K26").Value = 230.8

# Restore selection to match the saved view state
$ws.Range("J10").Select()
